{"js": "// 1) Insert a new \"Abstract\" styled paragraph right after the Title paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet titlePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].style === \"Title\") {\n    titlePara = paragraphs.items[i];\n    break;\n  }\n}\nif (!titlePara) {\n  titlePara = paragraphs.items[0];\n}\n\nconst abstractText =\n  \"my abstract Quarto enables you to weave together content and executable \" +\n  \"code into a finished document. To learn more about Quarto see https://quarto.org.\";\nconst abstractPara = titlePara.insertParagraph(abstractText, \"After\");\nabstractPara.style = \"Abstract\";\nawait context.sync();\n\n// 2) Rename the pander function call from pandoc.footnote.return to pandoc.footnote.\nconst fnNameHits = context.document.body.search(\"pandoc.footnote.return\", { matchCase: true });\nfnNameHits.load(\"items\");\nawait context.sync();\nif (fnNameHits.items.length > 0) {\n  fnNameHits.items[0].insertText(\"pandoc.footnote\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Strip the manual \"[1] \u201c \u2026 \u201d\" wrapper text around the footnote reference,\n//    leaving only the footnoteReference run in its paragraph.\nconst refParagraphs = context.document.body.paragraphs;\nrefParagraphs.load(\"items/text\");\nawait context.sync();\n\nlet footnoteWrapperPara = null;\nfor (let i = 0; i < refParagraphs.items.length; i++) {\n  if (refParagraphs.items[i].text.indexOf(\"[1] \\u201c\") === 0) {\n    footnoteWrapperPara = refParagraphs.items[i];\n    break;\n  }\n}\n\nif (footnoteWrapperPara) {\n  const footnotePrefixHits = footnoteWrapperPara.search(\"[1] \\u201c\", { matchCase: true });\n  footnotePrefixHits.load(\"items\");\n  await context.sync();\n  if (footnotePrefixHits.items.length > 0) {\n    footnotePrefixHits.items[0].insertText(\"\", \"Replace\");\n    await context.sync();\n  }\n\n  const footnoteSuffixHits = footnoteWrapperPara.search(\"\\u201d\", { matchCase: true });\n  footnoteSuffixHits.load(\"items\");\n  await context.sync();\n  if (footnoteSuffixHits.items.length > 0) {\n    footnoteSuffixHits.items[0].insertText(\"\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Insert a new \"Abstract\" styled paragraph right after the Title paragraph.\n$titlePara = $null\n$titleIndex = 0\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Style.NameLocal -eq \"Title\") {\n        $titlePara = $p\n        $titleIndex = $i\n        break\n    }\n}\nif ($titlePara -eq $null) {\n    $titlePara = $d.Paragraphs(1)\n    $titleIndex = 1\n}\n$titlePara.Range.InsertParagraphAfter()\n$abstractPara = $d.Paragraphs($titleIndex + 1)\n$abstractPara.Range.Text = \"my abstract Quarto enables you to weave together content and executable code into a finished document. To learn more about Quarto see https://quarto.org.\"\n$abstractPara.Style = \"Abstract\"\n\n# 2) Rename the pander function call from pandoc.footnote.return to pandoc.footnote.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"pandoc.footnote.return\"\n$find.Replacement.Text = \"pandoc.footnote\"\n$find.Execute($null, $true, $true, $false, $null, $null, $true, $null, $null, $find.Replacement.Text, 2) | Out-Null\n\n# 3) Strip the manual \"[1] \" ... \"\" wrapper text around the footnote reference,\n#    leaving only the footnoteReference run in its paragraph.\n$footnotePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains([char]0x201C)) {\n        $footnotePara = $p\n    }\n}\nif ($footnotePara -ne $null) {\n    $r = $footnotePara.Range\n    # Remove the trailing curly close-quote (right before the paragraph mark).\n    $suffixRange = $d.Range($r.End - 2, $r.End - 1)\n    $suffixRange.Text = \"\"\n    # Remove the leading \"[1] \" + curly open-quote.\n    $prefixRange = $d.Range($r.Start, $r.Start + 5)\n    $prefixRange.Text = \"\"\n}\n"}
